$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the shared-string text that appears in cell D1 ("LEA_LK" -> "TDK_MDT")
$ws.Range("D1").Value = "TDK_MDT"

# Delete the blank spacer column E (shifts F..N left by one column)
$ws.Columns.Item(5).Delete()

# Keep the previously-selected cell in sync with the column shift (was N7, now M7)
$ws.Range("M7").Select()
